# Auto-generated Excel COM-interop script
# Refreshes cached Universalis market-price / leve-profit columns (H:N) on several
# sheets, matching a scheduled market-data runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5175535
$ws.Range("I17").Value = 666
$ws.Range("J17").Value = 5606774
$ws.Range("K17").Value = 1998
$ws.Range("L17").Value = 16820322
$ws.Range("M17").Value = -1830
$ws.Range("N17").Value = -16820658
$ws.Range("H38").Value = 2077.2856
$ws.Range("I38").Value = 185.5
$ws.Range("J38").Value = 5860.857
$ws.Range("K38").Value = 556.5
$ws.Range("L38").Value = 17582.571
$ws.Range("M38").Value = -184.5
$ws.Range("N38").Value = -18326.571
$ws.Range("H40").Value = 2615.6667
$ws.Range("I40").Value = 1732.3334
$ws.Range("K40").Value = 1732.3334
$ws.Range("M40").Value = -1557.3334
$ws.Range("H62").Value = 7828
$ws.Range("I62").Value = 6959.2
$ws.Range("K62").Value = 6959.2
$ws.Range("M62").Value = -6335.2
$ws.Range("H65").Value = 7828
$ws.Range("I65").Value = 6959.2
$ws.Range("K65").Value = 34796
$ws.Range("M65").Value = -31676
$ws.Range("H74").Value = 5760.6924
$ws.Range("J74").Value = 5950
$ws.Range("L74").Value = 5950
$ws.Range("N74").Value = -7822
$ws.Range("H77").Value = 5760.6924
$ws.Range("J77").Value = 5950
$ws.Range("L77").Value = 29750
$ws.Range("N77").Value = -39110
$ws.Range("H100").Value = 2750.5
$ws.Range("I100").Value = 1376
$ws.Range("K100").Value = 1376
$ws.Range("M100").Value = -835
$ws.Range("H116").Value = 7061.6113
$ws.Range("I116").Value = 6575.125
$ws.Range("K116").Value = 6575.125
$ws.Range("M116").Value = -3133.125
$ws.Range("H137").Value = 40819812
$ws.Range("J137").Value = 50003770
$ws.Range("L137").Value = 150011310
$ws.Range("N137").Value = -150016410
$ws.Range("J138").Value = 9532913
$ws.Range("L138").Value = 28598739
$ws.Range("N138").Value = -28609019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15392178
$ws.Range("I32").Value = 16672393
$ws.Range("K32").Value = 16672393
$ws.Range("M32").Value = -16672106
$ws.Range("H74").Value = 200225820
$ws.Range("I74").Value = 250281780
$ws.Range("K74").Value = 250281780
$ws.Range("M74").Value = -250280906
$ws.Range("H77").Value = 200225820
$ws.Range("I77").Value = 250281780
$ws.Range("K77").Value = 1251408900
$ws.Range("M77").Value = -1251404532
$ws.Range("H132").Value = 25643792
$ws.Range("J132").Value = 111113270
$ws.Range("L132").Value = 333339810
$ws.Range("N132").Value = -333344870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 635.1429000000001
$ws.Range("I11").Value = 552.6667
$ws.Range("J11").Value = 697
$ws.Range("K11").Value = 552.6667
$ws.Range("L11").Value = 697
$ws.Range("M11").Value = -412.6667
$ws.Range("N11").Value = -977
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H86").Value = 22947.273
$ws.Range("I86").Value = 3689
$ws.Range("K86").Value = 3689
$ws.Range("M86").Value = -2566
$ws.Range("H89").Value = 22947.273
$ws.Range("I89").Value = 3689
$ws.Range("K89").Value = 18445
$ws.Range("M89").Value = -12829
$ws.Range("H97").Value = 18225
$ws.Range("I97").Value = 12258.429
$ws.Range("J97").Value = 59991
$ws.Range("K97").Value = 12258.429
$ws.Range("L97").Value = 59991
$ws.Range("M97").Value = -11267.429
$ws.Range("N97").Value = -61973
$ws.Range("H102").Value = 39939.5
$ws.Range("J102").Value = 69880
$ws.Range("L102").Value = 69880
$ws.Range("N102").Value = -76370
$ws.Range("H107").Value = 5976.3335
$ws.Range("I107").Value = 3599.1428
$ws.Range("K107").Value = 3599.1428
$ws.Range("M107").Value = -1679.1428
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2826
$ws.Range("I16").Value = 3031.5715
$ws.Range("J16").Value = 2106.5
$ws.Range("K16").Value = 3031.5715
$ws.Range("L16").Value = 2106.5
$ws.Range("M16").Value = -2744.5715
$ws.Range("N16").Value = -2680.5
$ws.Range("H86").Value = 3854.4546
$ws.Range("I86").Value = 2342.8572
$ws.Range("K86").Value = 2342.8572
$ws.Range("M86").Value = -1219.8572
$ws.Range("H89").Value = 3854.4546
$ws.Range("I89").Value = 2342.8572
$ws.Range("K89").Value = 11714.286
$ws.Range("M89").Value = -6098.286
$ws.Range("H99").Value = 12153.2
$ws.Range("J99").Value = 6248.5
$ws.Range("L99").Value = 6248.5
$ws.Range("N99").Value = -9244.5
$ws.Range("H104").Value = 49991
$ws.Range("J104").Value = 49991
$ws.Range("L104").Value = 49991
$ws.Range("N104").Value = -55233
$ws.Range("H113").Value = 2826
$ws.Range("I113").Value = 3031.5715
$ws.Range("J113").Value = 2106.5
$ws.Range("K113").Value = 3031.5715
$ws.Range("L113").Value = 2106.5
$ws.Range("M113").Value = -861.5715
$ws.Range("N113").Value = -6446.5
$ws.Range("H126").Value = 12153.2
$ws.Range("J126").Value = 6248.5
$ws.Range("L126").Value = 18745.5
$ws.Range("N126").Value = -23685.5
$ws.Range("H131").Value = 45000
$ws.Range("I131").Value = 45000
$ws.Range("K131").Value = 45000
$ws.Range("M131").Value = -39960
$ws.Range("H141").Value = 327150.9
$ws.Range("J141").Value = 350301
$ws.Range("L141").Value = 350301
$ws.Range("N141").Value = -360661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 673.0909
$ws.Range("J23").Value = 700.5
$ws.Range("L23").Value = 2101.5
$ws.Range("N23").Value = -2571.5
$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16080
$ws.Range("H82").Value = 7506
$ws.Range("I82").Value = 5013
$ws.Range("K82").Value = 15039
$ws.Range("M82").Value = -14633
$ws.Range("H85").Value = 7506
$ws.Range("I85").Value = 5013
$ws.Range("K85").Value = 15039
$ws.Range("M85").Value = -13635
$ws.Range("H87").Value = 3633
$ws.Range("I87").Value = 3633
$ws.Range("K87").Value = 10899
$ws.Range("M87").Value = -9651
$ws.Range("H90").Value = 3633
$ws.Range("I90").Value = 3633
$ws.Range("K90").Value = 32697
$ws.Range("M90").Value = -26457
$ws.Range("H128").Value = 116141.25
$ws.Range("I128").Value = 116141.25
$ws.Range("K128").Value = 348423.75
$ws.Range("M128").Value = -343443.75
$ws.Range("H136").Value = 5316.6665
$ws.Range("I136").Value = 10030
$ws.Range("J136").Value = 2960
$ws.Range("K136").Value = 30090
$ws.Range("L136").Value = 8880
$ws.Range("M136").Value = -24990
$ws.Range("N136").Value = -19080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5281.2104
$ws.Range("I70").Value = 4837.067
$ws.Range("J70").Value = 6946.75
$ws.Range("K70").Value = 4837.067
$ws.Range("L70").Value = 6946.75
$ws.Range("M70").Value = -4567.067
$ws.Range("N70").Value = -7486.75
$ws.Range("H73").Value = 5281.2104
$ws.Range("I73").Value = 4837.067
$ws.Range("J73").Value = 6946.75
$ws.Range("K73").Value = 4837.067
$ws.Range("L73").Value = 6946.75
$ws.Range("M73").Value = -3901.067
$ws.Range("N73").Value = -8818.75
$ws.Range("H102").Value = 2273.7368
$ws.Range("I102").Value = 2205.9375
$ws.Range("K102").Value = 2205.9375
$ws.Range("M102").Value = -583.9375
$ws.Range("H122").Value = 3166.6667
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -15400
$ws.Range("H126").Value = 5178045.5
$ws.Range("I126").Value = 2862838.8
$ws.Range("J126").Value = 8701186
$ws.Range("K126").Value = 8588516.399999999
$ws.Range("L126").Value = 26103558
$ws.Range("M126").Value = -8586046.399999999
$ws.Range("N126").Value = -26108498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2115.625
$ws.Range("I16").Value = 2115.625
$ws.Range("K16").Value = 2115.625
$ws.Range("M16").Value = -1945.625
$ws.Range("H46").Value = 1767.3
$ws.Range("J46").Value = 3378
$ws.Range("L46").Value = 3378
$ws.Range("N46").Value = -3754
$ws.Range("H93").Value = 5500
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 5500
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 5500
$ws.Range("N93").Value = -7996
$ws.Range("H112").Value = 53846.25
$ws.Range("J112").Value = 53846.25
$ws.Range("L112").Value = 53846.25
$ws.Range("N112").Value = -56800.25
$ws.Range("H131").Value = 54202
$ws.Range("J131").Value = 54202
$ws.Range("L131").Value = 54202
$ws.Range("N131").Value = -64282
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47670972
$ws.Range("I122").Value = 58886390
$ws.Range("K122").Value = 176659170
$ws.Range("M122").Value = -176656720

Write-Host "Applied 233 value updates and 2 cleared cells across 8 sheets."
